$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row cells: "q" -> "question", "a" -> "answer"
$ws.Range("A1").Value = "question"
$ws.Range("B1").Value = "answer"

# Update the selected cell/range to match the new view state
$ws.Range("C4").Select()
